$wb = $excel.ActiveWorkbook
$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1. Header / title font re-colour (bold white text) on both dashboards.
#    The old 14pt title font and the plain bold header font are effectively
#    merged into a single bold + white font used by both the title row and
#    the table header row.
# ---------------------------------------------------------------------------
$wsTraining.Range("A1").Font.Bold = $true
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsTraining.Range("A2:K2").Font.Bold = $true
$wsTraining.Range("A2:K2").Font.Color = 16777215

$wsExam.Range("A1").Font.Bold = $true
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Bold = $true
$wsExam.Range("A2:G2").Font.Color = 16777215

# ---------------------------------------------------------------------------
# 2. Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
#    columns - 8 more days have passed since the report was generated.
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3 = 269; 4 = 308; 5 = 267; 6 = 363; 7 = 244; 8 = 336;
    9 = 251; 10 = 265; 11 = 266; 12 = 245; 13 = 323; 14 = 349;
    15 = 260; 16 = 413; 17 = 413; 18 = -23; 19 = -83; 20 = -106;
    21 = -34; 22 = -34; 23 = 155; 24 = 268
}

for ($r = 3; $r -le 24; $r++) {
    $wsTraining.Cells.Item($r, 8).Value = $periodToExpire[$r]
    # Use the leading-apostrophe text prefix so Excel stores the new date as
    # plain text (matching the original inline string) instead of silently
    # re-interpreting "16-Sep-2025" as a real date value.
    $wsTraining.Cells.Item($r, 9).Formula = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3. Exam Dashboard: narrow the COMMENTS column and normalise its text.
# ---------------------------------------------------------------------------
$wsExam.Columns.Item(5).ColumnWidth = 14.1666666667

for ($r = 3; $r -le 11; $r++) {
    $wsExam.Cells.Item($r, 5).Value = "date is valid"
}
